$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Typography")
$ws2 = $wb.Worksheets.Item("Translation")

$ws1.Range("F4").Value = "?"
$ws1.Range("G4").Value = "0x20-0xFE"
$ws1.Range("H4").Value = "0x20-0xFE"

$ws1.Range("F5").Value = "?"
$ws1.Range("G5").Value = "0x20-0xFE"
$ws1.Range("H5").Value = "0x20-0xFE"

$ws1.Range("F6").Value = "?"
$ws1.Range("G6").Value = "0-9,A-Z,a-z"
$ws1.Range("H6").Value = "0-9,A-Z,a-z"

$ws2.Range("E4").Value = "change text"

$ws2.Range("B7").Value = "SingleUseId4"
$ws2.Range("C7").Value = "Default"
$ws2.Range("D7").Value = "Left"
$ws2.Range("E7").Value = "test text"
$ws2.Range("F7").Value = "LTR"

$ws2.Range("B8").Value = "SingleUseId5"
$ws2.Range("C8").Value = "Default"
$ws2.Range("D8").Value = "Left"
$ws2.Range("E8").Value = "UART TEST"
$ws2.Range("F8").Value = "LTR"

$ws2.Range("B9").Value = "SingleUseId6"
$ws2.Range("C9").Value = "Default"
$ws2.Range("D9").Value = "Center"
$ws2.Range("E9").Value = "UART INIT"
$ws2.Range("F9").Value = "LTR"

$ws2.Range("B10").Value = "SingleUseId7"
$ws2.Range("C10").Value = "Default"
$ws2.Range("D10").Value = "Center"
$ws2.Range("E10").Value = "SEND"
$ws2.Range("F10").Value = "LTR"

$ws2.Range("B7:F10").Style = "Normal"
